$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# NOTIFY_SEND_SETTING sheet: add 3 new rows (19-21) for the scheduled job
# status notifications (site / organization / project level).
# ---------------------------------------------------------------------------
$wsSetting = $wb.Worksheets.Item("NOTIFY_SEND_SETTING")

# Row 19 - jobStatusSite
$wsSetting.Range("E19").Value = "jobStatusSite"
$wsSetting.Range("F19").Value = "jobStatusSite"
$wsSetting.Range("G19").Value = "平台任务状态通知"
$wsSetting.Range("H19").Value = "定时任务状态发生变化，给相关用户发送通知。"
$wsSetting.Range("I19").Value = "site"
$wsSetting.Range("J19").Value = 0
$wsSetting.Range("K19").Value = 1
$wsSetting.Range("L19").Value = 1
$wsSetting.Range("M19").Value = 1

# Row 20 - jobStatusOrganization
$wsSetting.Range("E20").Value = "jobStatusOrganization"
$wsSetting.Range("F20").Value = "jobStatusOrganization"
$wsSetting.Range("G20").Value = "组织任务状态通知"
$wsSetting.Range("H20").Value = "定时任务状态发生变化，给相关用户发送通知。"
$wsSetting.Range("I20").Value = "organization"
$wsSetting.Range("J20").Value = 0
$wsSetting.Range("K20").Value = 1
$wsSetting.Range("L20").Value = 1
$wsSetting.Range("M20").Value = 1

# Row 21 - jobStatusProject
$wsSetting.Range("E21").Value = "jobStatusProject"
$wsSetting.Range("F21").Value = "jobStatusProject"
$wsSetting.Range("G21").Value = "项目任务状态通知"
$wsSetting.Range("H21").Value = "定时任务状态发生变化，给相关用户发送通知。"
$wsSetting.Range("I21").Value = "project"
$wsSetting.Range("J21").Value = 0
$wsSetting.Range("K21").Value = 1
$wsSetting.Range("L21").Value = 1
$wsSetting.Range("M21").Value = 1

# ---------------------------------------------------------------------------
# NOTIFY_TEMPLATE sheet: add 3 new rows (19-21) with the PM message templates
# for the same 3 new notification business types.
# ---------------------------------------------------------------------------
$wsTemplate = $wb.Worksheets.Item("NOTIFY_TEMPLATE")

# Row 19 - jobStatusSite-preset
$wsTemplate.Range("E19").Value = 12
$wsTemplate.Range("F19").Value = "jobStatusSite-preset"
$wsTemplate.Range("G19").Value = "任务状态变更通知"
$wsTemplate.Range("H19").Value = "pm"
$wsTemplate.Range("I19").Value = 1
$wsTemplate.Range("J19").Value = "jobStatusSite"
$wsTemplate.Range("L19").Value = '${jobName}状态发生变更'
$wsTemplate.Range("M19").Value = '您好，${userName}。${jobName}任务已${jobStatus}，请注意查看。'

# Row 20 - jobStatusOrganization-preset
$wsTemplate.Range("E20").Value = 13
$wsTemplate.Range("F20").Value = "jobStatusOrganization-preset"
$wsTemplate.Range("G20").Value = "任务状态变更通知"
$wsTemplate.Range("H20").Value = "pm"
$wsTemplate.Range("I20").Value = 1
$wsTemplate.Range("J20").Value = "jobStatusOrganization"
$wsTemplate.Range("L20").Value = '${jobName}状态发生变更'
$wsTemplate.Range("M20").Value = '您好，${userName}。${jobName}任务已${jobStatus}，请注意查看。'

# Row 21 - jobStatusProject-preset
$wsTemplate.Range("E21").Value = 14
$wsTemplate.Range("F21").Value = "jobStatusProject-preset"
$wsTemplate.Range("G21").Value = "任务状态变更通知"
$wsTemplate.Range("H21").Value = "pm"
$wsTemplate.Range("I21").Value = 1
$wsTemplate.Range("J21").Value = "jobStatusProject"
$wsTemplate.Range("L21").Value = '${jobName}状态发生变更'
$wsTemplate.Range("M21").Value = '您好，${userName}。${jobName}任务已${jobStatus}，请注意查看。'
